# Final version - Bugs corrected
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "longueur_prolongement_cotes (mm)" row (row 10) entirely -
# rows below shift up, so epaisseur_plastique_fin/largeur_prolongement_cotes/
# hauteur_accroche move from rows 11-13 to rows 10-12.
$ws.Rows.Item(10).Delete()

# Reselect a cell below the used range, matching the saved selection state.
$ws.Range("B17").Select()

# Set the page setup (paper size + orientation), which is what stamps a
# pageSetup element + printer-settings relationship onto the sheet part.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
